$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Stage" column header in E1. Start from D1's header formatting
# (bold white font on the blue fill) by copy/paste-format so we reuse the
# existing font/fill entries instead of minting new ones, then strip the
# outer border that D1 has (the new header style has no border).
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").Value = "Stage"
$ws.Range("E1").Borders.LineStyle = -4142

# New column width for column E (closest value the host's pixel-quantised
# ColumnWidth setter can reproduce for the authored stored width of
# 29.26953125 characters).
$ws.Range("E1").ColumnWidth = 28.5

# Update selection to E2
$ws.Range("E2").Select()
